$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> column -> new text value (as string, since cells are text/inlineStr)
$updates = @{
    2 = @{ 'D'='305.96'; 'E'='6.76%'; 'G'='2' }
    3 = @{ 'D'='32.02'; 'E'='9.20%'; 'G'='2' }
    4 = @{ 'D'='5.307'; 'E'='4.36%'; 'G'='2' }
    5 = @{ 'D'='0.07415'; 'E'='11.42%'; 'G'='2' }
    6 = @{ 'D'='7.782'; 'E'='6.09%'; 'G'='2' }
    7 = @{ 'D'='3.684'; 'E'='8.16%'; 'G'='2' }
    8 = @{ 'D'='1.460'; 'E'='8.76%'; 'G'='2' }
    9 = @{ 'D'='0.9126'; 'E'='-1.23%'; 'G'='2' }
    10 = @{ 'D'='0.01630'; 'E'='2,425.56%'; 'G'='2' }
    11 = @{ 'D'='0.1687'; 'E'='7.81%'; 'G'='2' }
    12 = @{ 'D'='0.07663'; 'E'='21.17%'; 'G'='2' }
    13 = @{ 'D'='0.08030'; 'E'='5.35%'; 'G'='2' }
    14 = @{ 'D'='0.03029'; 'E'='4.46%'; 'G'='2' }
    15 = @{ 'D'='0.09807'; 'E'='9.14%'; 'G'='2' }
    16 = @{ 'D'='0.001515'; 'E'='-5.25%'; 'G'='2' }
    17 = @{ 'D'='0.006359'; 'E'='1.65%'; 'G'='2' }
    18 = @{ 'D'='3.506'; 'E'='1.40%'; 'G'='2' }
    19 = @{ 'D'='2.247'; 'E'='0.79%'; 'G'='2' }
    20 = @{ 'D'='0.3269'; 'E'='1.74%'; 'G'='2' }
    21 = @{ 'D'='0.1296'; 'E'='-1.03%'; 'G'='2' }
    22 = @{ 'D'='4.241'; 'E'='4.51%'; 'G'='2' }
    23 = @{ 'D'='0.1611'; 'E'='3.82%'; 'G'='2' }
    24 = @{ 'D'='0.04532'; 'E'='1.09%'; 'G'='2' }
    25 = @{ 'D'='0.001218'; 'E'='2.25%'; 'G'='2' }
    26 = @{ 'D'='0.004511'; 'E'='9.57%'; 'G'='2' }
    27 = @{ 'D'='0.0001166'; 'E'='-6.65%'; 'G'='2' }
    28 = @{ 'D'='0.0001769'; 'E'='9.33%'; 'G'='2' }
    29 = @{ 'G'='2' }
    30 = @{ 'G'='2' }
    31 = @{ 'G'='2' }
    32 = @{ 'G'='2' }
    33 = @{ 'G'='2' }
    34 = @{ 'G'='2' }
    35 = @{ 'G'='2' }
    36 = @{ 'G'='2' }
    37 = @{ 'G'='2' }
    38 = @{ 'G'='2' }
    39 = @{ 'G'='2' }
    40 = @{ 'D'='0.04517'; 'E'='8.02%'; 'G'='2' }
    41 = @{ 'D'='0.007100'; 'E'='5.06%'; 'G'='2' }
    42 = @{ 'D'='0.1360'; 'E'='9.85%'; 'G'='2' }
    43 = @{ 'D'='0.002172'; 'E'='9.80%'; 'G'='2' }
    44 = @{ 'D'='0.01345'; 'E'='7.59%'; 'G'='2' }
    45 = @{ 'D'='0.00005921'; 'E'='5.78%'; 'G'='2' }
    46 = @{ 'D'='1.896'; 'E'='-3.66%'; 'G'='2' }
    47 = @{ 'D'='0.01295'; 'E'='-0.90%'; 'G'='2' }
    48 = @{ 'G'='2' }
    49 = @{ 'G'='2' }
    50 = @{ 'G'='2' }
    51 = @{ 'G'='2' }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
    }
}
